$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("list 1")
$ws1.Range("B1:C1").EntireColumn.Insert()
$ws1.Range("B2:C2").Font.Name = "Calibri"
$ws1.Range("B11:C11").Font.Name = "Calibri"
$ws1.Range("B1").Value = "en_comments"
$ws1.Range("C1").Value = "de_comments"
$ws1.Range("B2").Value = "Comment for the list."
$ws1.Range("C2").Value = "Kommentar für die Liste"
$ws1.Range("B11").Value = "Comment for the third node."
$ws1.Range("C11").Value = "Kommentar für den dritten Knoten."

$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("B1:F1").EntireColumn.Insert()
$ws2.Range("B1").Value = "en_comments"
$ws2.Range("C1").Value = "de_comments"
$ws2.Range("D1").Value = "fr_comments"
$ws2.Range("E1").Value = "it_comments"
$ws2.Range("F1").Value = "rm_comments"
